# Refresh cryptos worksheet with the latest scrape (GitHub Actions run).
# D = Price text, E = 1h volume-change text; a handful of rows also
# re-ranked (coin name/link/price/volume all move together).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  # row 2: Bitcoin
    @("D2", '30.747.70'),
    @("E2", '  +2.48%  '),
  # row 3: Ethereum
    @("D3", '1.895.43'),
    @("E3", '  +0.77%  '),
  # row 4: TetherUSD
    @("D4", '''1.001'),
    @("E4", '  +0.24%  '),
  # row 5: BNB
    @("D5", '''247.70'),
    @("E5", '  +1.91%  '),
  # row 6: USDC
    @("D6", '''0.9999'),
    @("E6", '  +0.19%  '),
  # row 7: XRP
    @("D7", '''0.4943'),
    @("E7", '  -0.34%  '),
  # row 8: Cardano
    @("D8", '''0.2968'),
    @("E8", '  +1.61%  '),
  # row 9: Dogecoin
    @("D9", '''0.06829'),
    @("E9", '  +2.86%  '),
  # row 10: WrappedEther
    @("D10", '1.894.85'),
    @("E10", '  +0.77%  '),
  # row 11: Solana
    @("E11", '  +3.36%  '),
  # row 12: Litecoin
    @("D12", '''92.59'),
    @("E12", '  +7.07%  '),
  # row 13: TRON
    @("D13", '''0.07276'),
    @("E13", '  +0.54%  '),
  # row 14: Polkadot
    @("D14", '''5.114'),
    @("E14", '  +4.82%  '),
  # row 15: Polygon
    @("D15", '''0.6816'),
    @("E15", '  +1.84%  '),
  # row 16: WrappedBTC
    @("D16", '30.725.98'),
    @("E16", '  +2.53%  '),
  # row 17: ShibaInu
    @("D17", '''0.000008010'),
    @("E17", '  +1.23%  '),
  # row 18: Avalanche
    @("E18", '  +4.30%  '),
  # row 19: Dai
    @("D19", '''0.9997'),
    @("E19", '  +0.14%  '),
  # row 20: WrappedliquidstakedEther2.0
    @("D20", '2.139.27'),
    @("E20", '  +0.72%  '),
  # row 21: BinanceUSD
    @("D21", '''1.000'),
    @("E21", '  +0.31%  '),
  # row 22: Uniswap
    @("D22", '''4.865'),
    @("E22", '  +1.96%  '),
  # row 23: BitcoinCash
    @("D23", '''193.87'),
    @("E23", '  +34.83%  '),
  # row 24: Chainlink
    @("D24", '''6.101'),
    @("E24", '  +7.42%  '),
  # row 25: Cosmos
    @("D25", '''9.471'),
    @("E25", '  +4.46%  '),
  # row 26: Monero
    @("D26", '''155.39'),
    @("E26", '  +3.64%  '),
  # row 27: EthereumClassic
    @("D27", '''19.28'),
    @("E27", '  +12.54%  '),
  # row 28: LidoDAOToken
    @("D28", '''1.930'),
    @("E28", '  +0.64%  '),
  # row 29: Toncoin
    @("E29", '  +0.71%  '),
  # row 30: InternetComputer(DFINITY)
    @("D30", '''4.369'),
    @("E30", '  +4.24%  '),
  # row 31: Stellar
    @("E31", '  +3.19%  '),
  # row 32: Filecoin
    @("D32", '''4.053'),
    @("E32", '  +2.74%  '),
  # row 33: Hedera
    @("E33", '  +2.76%  '),
  # row 34: ImmutableX
    @("D34", '''0.7465'),
    @("E34", '  +4.67%  '),
  # row 35: ARBITRUM
    @("E35", '  +1.45%  '),
  # row 36: HuobiToken
    @("D36", '''2.734'),
    @("E36", '  +2.56%  '),
  # row 37: VeChain
    @("E37", '  +4.62%  '),
  # row 38: MXToken
    @("D38", '''2.687'),
    @("E38", '  -0.16%  '),
  # row 39: RenderToken
    @("D39", '''2.173'),
    @("E39", '  -0.24%  '),
  # row 40: TrustWalletToken
    @("D40", '''0.9423'),
    @("E40", '  +1.05%  '),
  # row 41: TheSandbox
    @("D41", '''0.4453'),
    @("E41", '  +4.84%  '),
  # row 42: Quant
    @("D42", '''106.55'),
    @("E42", '  +4.01%  '),
  # row 43: FraxShare
    @("D43", '''5.772'),
    @("E43", '  -0.60%  '),
  # row 44: PaxDollar
    @("D44", '''0.9999'),
    @("E44", '  +0.19%  '),
  # row 45: Aptos
    @("D45", '''7.720'),
    @("E45", '  +3.17%  '),
  # row 46: Algorand
    @("D46", '''0.1348'),
    @("E46", '  +6.88%  '),
  # row 47: Cronos
    @("E47", '  +3.75%  '),
  # row 48: NEARProtocol/EnergySwap (swapped)
    @("B48", 'EnergySwap'),
    @("C48", 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'),
    @("D48", '''8.766'),
    @("E48", '  +5.97%  '),
  # row 49: EnergySwap/NEARProtocol (swapped)
    @("B49", 'NEARProtocol'),
    @("C49", 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'),
    @("D49", '''1.436'),
    @("E49", '  +7.55%  '),
  # row 50: Decentraland
    @("D50", '''0.3958'),
    @("E50", '  +5.02%  '),
  # row 51: Elrond
    @("D51", '''33.69'),
    @("E51", '  +3.85%  ')
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
